$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A83").Value = "GRT-USD"
